# Scheduled runner update: refresh leve-profit calculations (currentAveragePrice,
# NQ/HQ prices and profits) across all 8 job sheets with newly recomputed
# market-board figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 516.2857
$ws.Range("I2").Value = 499.6
$ws.Range("K2").Value = 499.6
$ws.Range("M2").Value = -386.6

$ws.Range("H38").Value = 12232.923
$ws.Range("I38").Value = 16718.715
$ws.Range("J38").Value = 6999.5
$ws.Range("K38").Value = 50156.145
$ws.Range("L38").Value = 20998.5
$ws.Range("M38").Value = -49784.145
$ws.Range("N38").Value = -21742.5

$ws.Range("H58").Value = 19234824
$ws.Range("J58").Value = 23250
$ws.Range("L58").Value = 69750
$ws.Range("N58").Value = -70050

$ws.Range("H132").Value = 1789.475
$ws.Range("I132").Value = 1385.1428
$ws.Range("K132").Value = 4155.428400000001
$ws.Range("M132").Value = -1625.428400000001

$ws.Range("H133").Value = 64999.2
$ws.Range("J133").Value = 64999.2
$ws.Range("L133").Value = 64999.2
$ws.Range("N133").Value = -75119.2

$ws.Range("H139").Value = 70225.32000000001
$ws.Range("J139").Value = 70225.32000000001
$ws.Range("L139").Value = 70225.32000000001
$ws.Range("N139").Value = -80505.32000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 933
$ws.Range("I12").Value = 799
$ws.Range("K12").Value = 799
$ws.Range("M12").Value = -626

$ws.Range("H74").Value = 1627.6522
$ws.Range("I74").Value = 1616
$ws.Range("J74").Value = 1750
$ws.Range("K74").Value = 1616
$ws.Range("L74").Value = 1750
$ws.Range("M74").Value = -742
$ws.Range("N74").Value = -3498

$ws.Range("H77").Value = 1627.6522
$ws.Range("I77").Value = 1616
$ws.Range("J77").Value = 1750
$ws.Range("K77").Value = 8080
$ws.Range("L77").Value = 8750
$ws.Range("M77").Value = -3712
$ws.Range("N77").Value = -17486

$ws.Range("H122").Value = 12348736

$ws.Range("H132").Value = 2260.5293
$ws.Range("I132").Value = 2260.5293
$ws.Range("K132").Value = 6781.5879
$ws.Range("M132").Value = -4251.5879

$ws.Range("H133").Value = 99999.836
$ws.Range("J133").Value = 99999.836
$ws.Range("L133").Value = 99999.836
$ws.Range("N133").Value = -105059.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 71442430
$ws.Range("I86").Value = 15500
$ws.Range("J86").Value = 100013200
$ws.Range("K86").Value = 15500
$ws.Range("L86").Value = 100013200
$ws.Range("M86").Value = -14377
$ws.Range("N86").Value = -100015446

$ws.Range("H89").Value = 71442430
$ws.Range("I89").Value = 15500
$ws.Range("J89").Value = 100013200
$ws.Range("K89").Value = 77500
$ws.Range("L89").Value = 500066000
$ws.Range("M89").Value = -71884
$ws.Range("N89").Value = -500077232

$ws.Range("H94").Value = 2464.6667
$ws.Range("I94").Value = 2259.6191
$ws.Range("J94").Value = 3900
$ws.Range("K94").Value = 2259.6191
$ws.Range("L94").Value = 3900
$ws.Range("M94").Value = -1808.6191
$ws.Range("N94").Value = -4802

$ws.Range("H97").Value = 12500
$ws.Range("I97").Value = 12500
$ws.Range("K97").Value = 12500
$ws.Range("M97").Value = -11509

$ws.Range("H134").Value = 2694.9412
$ws.Range("I134").Value = 2787.5
$ws.Range("K134").Value = 8362.5
$ws.Range("M134").Value = -5827.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2431.3
$ws.Range("I132").Value = 2439.8235
$ws.Range("K132").Value = 7319.470499999999
$ws.Range("M132").Value = -4789.470499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 177117.77
$ws.Range("I121").Value = 501.8
$ws.Range("J121").Value = 429426.28
$ws.Range("K121").Value = 1505.4
$ws.Range("L121").Value = 1288278.84
$ws.Range("M121").Value = -195.4000000000001
$ws.Range("N121").Value = -1290898.84

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2355
$ws.Range("I132").Value = 2445.125
$ws.Range("K132").Value = 7335.375
$ws.Range("M132").Value = -4805.375

$ws.Range("H139").Value = 97583.39999999999
$ws.Range("J139").Value = 97583.39999999999
$ws.Range("L139").Value = 97583.39999999999
$ws.Range("N139").Value = -107863.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3110.1428
$ws.Range("I7").Value = 3110.1428
$ws.Range("K7").Value = 3110.1428
$ws.Range("M7").Value = -2998.1428

$ws.Range("H22").Value = 3072.0908
$ws.Range("I22").Value = 2406.4
$ws.Range("J22").Value = 3626.8333
$ws.Range("K22").Value = 2406.4
$ws.Range("L22").Value = 3626.8333
$ws.Range("M22").Value = -2111.4
$ws.Range("N22").Value = -4216.8333

$ws.Range("H27").Value = 3072.0908
$ws.Range("I27").Value = 2406.4
$ws.Range("J27").Value = 3626.8333
$ws.Range("K27").Value = 2406.4
$ws.Range("L27").Value = 3626.8333
$ws.Range("M27").Value = -2299.4
$ws.Range("N27").Value = -3840.8333

$ws.Range("H34").Value = 14024
$ws.Range("J34").Value = 14024
$ws.Range("L34").Value = 14024
$ws.Range("N34").Value = -14368

$ws.Range("H74").Value = 52000
$ws.Range("I74").Value = 52000
$ws.Range("K74").Value = 52000
$ws.Range("M74").Value = -51002

$ws.Range("H77").Value = 52000
$ws.Range("I77").Value = 52000
$ws.Range("K77").Value = 156000
$ws.Range("M77").Value = -151008

$ws.Range("H93").Value = 5895.2104
$ws.Range("J93").Value = 6800.6
$ws.Range("L93").Value = 6800.6
$ws.Range("N93").Value = -9296.6

$ws.Range("H100").Value = 7899.909
$ws.Range("I100").Value = 5479.8
$ws.Range("K100").Value = 5479.8
$ws.Range("M100").Value = -4938.8

$ws.Range("H126").Value = 3110.1428
$ws.Range("I126").Value = 3110.1428
$ws.Range("K126").Value = 9330.428400000001
$ws.Range("M126").Value = -6860.428400000001

$ws.Range("H131").Value = 73326
$ws.Range("J131").Value = 73326
$ws.Range("L131").Value = 73326
$ws.Range("N131").Value = -83406

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 40002400
$ws.Range("J5").Value = 40002400
$ws.Range("L5").Value = 40002400
$ws.Range("N5").Value = -40002624

$ws.Range("H75").Value = 50000
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 50000
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H99").Value = 50000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H126").Value = 2539.6
$ws.Range("I126").Value = 2623.25
$ws.Range("K126").Value = 7869.75
$ws.Range("M126").Value = -5399.75

$ws.Range("H132").Value = 7323.5654
$ws.Range("I132").Value = 7530.6665
$ws.Range("K132").Value = 22591.9995
$ws.Range("M132").Value = -20061.9995
